$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "OrderNo"
$ws.Range("B1").Value = "Product"
$ws.Range("C1").Value = "Fulfilment Store"
$ws.Range("D1").Value = "Total Price"

# D2 cell gets a currency number format (builtin format 6: "$"#,##0;[Red]-"$"#,##0)
$ws.Range("D2").NumberFormat = "$#,##0_);[Red]($#,##0)"

# Update selection to D7 as shown in diff
$ws.Range("D7").Select()
